$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.855.61'
$ws.Range("E2").Value = '  +1.70%  '
$ws.Range("D3").Value = '1.728.04'
$ws.Range("E3").Value = '  +0.25%  '
$ws.Range("D4").Value = '''0.9984'
$ws.Range("E4").Value = '  -0.18%  '
$ws.Range("D5").Value = '''241.73'
$ws.Range("E5").Value = '  -0.73%  '
$ws.Range("D6").Value = '''0.9985'
$ws.Range("E6").Value = '  -0.22%  '
$ws.Range("D7").Value = '''0.4895'
$ws.Range("E7").Value = '  -0.17%  '
$ws.Range("D8").Value = '''0.2595'
$ws.Range("E8").Value = '  -0.48%  '
$ws.Range("D9").Value = '''0.06218'
$ws.Range("E9").Value = '  +0.41%  '
$ws.Range("D10").Value = '1.732.43'
$ws.Range("E10").Value = '  +0.45%  '
$ws.Range("D11").Value = '''16.00'
$ws.Range("E11").Value = '  +3.28%  '
$ws.Range("D12").Value = '''0.06916'
$ws.Range("E12").Value = '  -1.20%  '
$ws.Range("D13").Value = '''0.6087'
$ws.Range("E13").Value = '  +1.61%  '
$ws.Range("D14").Value = '''4.491'
$ws.Range("E14").Value = '  -1.75%  '
$ws.Range("D15").Value = '''77.29'
$ws.Range("D16").Value = '''0.9988'
$ws.Range("D17").Value = '26.638.61'
$ws.Range("E17").Value = '  +0.87%  '
$ws.Range("D18").Value = '''0.9983'
$ws.Range("E18").Value = '  -0.20%  '
$ws.Range("D19").Value = '''0.000007179'
$ws.Range("E19").Value = '  +0.71%  '
$ws.Range("D20").Value = '''11.44'
$ws.Range("E20").Value = '  +0.93%  '
$ws.Range("D21").Value = '1.957.72'
$ws.Range("E21").Value = '  +0.63%  '
$ws.Range("D22").Value = '''4.431'
$ws.Range("E22").Value = '  -0.76%  '
$ws.Range("D23").Value = '''8.564'
$ws.Range("E23").Value = '  -0.22%  '
$ws.Range("D24").Value = '''5.120'
$ws.Range("E24").Value = '  -0.63%  '
$ws.Range("D25").Value = '''138.40'
$ws.Range("E25").Value = '  +0.75%  '
$ws.Range("D26").Value = '''15.32'
$ws.Range("E26").Value = '  +0.72%  '
$ws.Range("D27").Value = '''1.777'
$ws.Range("E27").Value = '  +4.72%  '
$ws.Range("D28").Value = '''1.383'
$ws.Range("E28").Value = '  -0.57%  '
$ws.Range("D29").Value = '''106.33'
$ws.Range("E29").Value = '  -0.48%  '
$ws.Range("D30").Value = '''3.949'
$ws.Range("E30").Value = '  +0.26%  '
$ws.Range("D31").Value = '''0.07990'
$ws.Range("E31").Value = '  +0.66%  '
$ws.Range("D32").Value = '''3.687'
$ws.Range("E32").Value = '  +0.43%  '
$ws.Range("D33").Value = '''0.04535'
$ws.Range("E33").Value = '  +0.03%  '
$ws.Range("B34").Value = 'Frax'
$ws.Range("C34").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D34").Value = '''0.9982'
$ws.Range("E34").Value = '  -0.14%  '
$ws.Range("B35").Value = 'HuobiToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D35").Value = '''2.611'
$ws.Range("E35").Value = '  +0.28%  '
$ws.Range("B36").Value = 'ARBITRUM'
$ws.Range("C36").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D36").Value = '''1.009'
$ws.Range("E36").Value = '  +1.50%  '
$ws.Range("B37").Value = 'ImmutableX'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D37").Value = '''0.6237'
$ws.Range("E37").Value = '  -0.35%  '
$ws.Range("B38").Value = 'TrustWalletToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D38").Value = '''0.9354'
$ws.Range("E38").Value = '  +2.43%  '
$ws.Range("B39").Value = 'RenderToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D39").Value = '''2.052'
$ws.Range("E39").Value = '  +4.82%  '
$ws.Range("B40").Value = 'MXToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D40").Value = '''2.449'
$ws.Range("E40").Value = '  +2.37%  '
$ws.Range("B41").Value = 'PaxDollar'
$ws.Range("C41").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D41").Value = '''0.9994'
$ws.Range("E41").Value = '  -0.16%  '
$ws.Range("B42").Value = 'VeChain'
$ws.Range("C42").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D42").Value = '''0.01503'
$ws.Range("E42").Value = '  +1.48%  '
$ws.Range("B43").Value = 'FraxShare'
$ws.Range("C43").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D43").Value = '''5.668'
$ws.Range("E43").Value = '  +4.38%  '
$ws.Range("B44").Value = 'Quant'
$ws.Range("C44").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D44").Value = '''99.52'
$ws.Range("E44").Value = '  -0.39%  '
$ws.Range("B45").Value = 'TheSandbox'
$ws.Range("C45").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D45").Value = '''0.3859'
$ws.Range("E45").Value = '  +0.72%  '
$ws.Range("B46").Value = 'Aptos'
$ws.Range("C46").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D46").Value = '''6.930'
$ws.Range("E46").Value = '  +3.53%  '
$ws.Range("B47").Value = 'Algorand'
$ws.Range("C47").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D47").Value = '''0.1158'
$ws.Range("E47").Value = '  +0.17%  '
$ws.Range("B48").Value = 'Cronos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D48").Value = '''0.05391'
$ws.Range("E48").Value = '  +0.53%  '
$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D49").Value = '''7.969'
$ws.Range("E49").Value = '  +3.53%  '
$ws.Range("B50").Value = 'Elrond'
$ws.Range("C50").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D50").Value = '''30.19'
$ws.Range("E50").Value = '  +0.34%  '
$ws.Range("B51").Value = 'NEARProtocol'
$ws.Range("C51").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D51").Value = '''1.239'
$ws.Range("E51").Value = '  +0.27%  '
